# Rename the input sheets to carry an "input_" prefix (the "component_names"
# sheet keeps its original name).
$wb = $excel.ActiveWorkbook

$wsStoich = $wb.Worksheets.Item(1)
$wsStoich.Name = "input_stoich_coefficients"

$wsConc = $wb.Worksheets.Item(2)
$wsConc.Name = "input_concentrations"

$wsK = $wb.Worksheets.Item(3)
$wsK.Name = "input_k_constants_log10"

# Update each sheet's remembered selection / active cell, and make the
# k_constants_log10 sheet the active tab (matching the new activeTab=2 /
# tabSelected on sheet3 in the saved workbook).
$wsStoich.Activate()
$wsStoich.Range("J34").Select()

$wsConc.Activate()
$wsConc.Range("I33").Select()

$wsK.Activate()
$wsK.Range("L34").Select()
